$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Status column updates -------------------------------------------------
# Row 6 (task "5. Search algorithm optimization"): was Done, is reopened ->
# drop the green "Done" colour back to the sheet's normal automatic colour.
$ws.Range("C6").Value = "Reopened to following optimization"
$ws.Range("C6").Font.ThemeColor = 1

# Row 4 (task "3. Rename ...") moves from "In Progress" to "Done" -> green text.
$ws.Range("C4").Value = "Done"
$ws.Range("C4").Font.Color = 5287936

# Row 8 (task "7. Remove head link ...") moves from "Open" to "Done" -> green
# text, and gets an assignee.
$ws.Range("C8").Value = "Done"
$ws.Range("C8").Font.Color = 5287936
$ws.Range("D8").Value = "Arthur"

# --- New task row (row 15 was the first blank row in the table) ------------
$ws.Range("A15").Value = "14. Fix bug with dissappearing tags "
$ws.Range("B15").Value = "Hight"
$ws.Range("C15").Value = "Open"

# --- Grow the table with 6 more blank rows before the closing border row ---
$ws.Rows.Item(21).Resize(6).Insert()
$ws.Range("A15:D20").Copy()
$ws.Range("A21:D26").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Selection moves along with the edits -----------------------------------
$ws.Range("E13").Select()
